$d = $word.ActiveDocument

# Replace the text of the first paragraph with the new 2025 wording.
$d.Content.Find.Execute("Hola, me llaman Romeo <3.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Necesito lograr mi libertad este 2025.", 2)

# Append a blank paragraph after the (now edited) first paragraph.
$r1 = $d.Content
$r1.Collapse(0)
$r1.InsertParagraphAfter()

# Append another paragraph after the blank one, which will hold the new text.
$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# Put the closing remark text into that last (currently empty) paragraph.
$r3 = $d.Content
$r3.Collapse(0)
$r3.InsertAfter("(Así sea parcialmente)")
